$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 37, shifting existing rows 37:84 down to 38:85.
$ws.Range("A37").EntireRow.Insert()

# Populate the newly inserted row 37 with the new weekly record.
$ws.Cells.Item(37, 1).Value = 11
$ws.Cells.Item(37, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(37, 3).Value = "Bíobío"
$ws.Cells.Item(37, 4).Value = 44763
$ws.Cells.Item(37, 5).Value = 8
$ws.Cells.Item(37, 6).Value = 100112012
$ws.Cells.Item(37, 7).Value = "Espinaca"
$ws.Cells.Item(37, 8).Value = "Sin especificar"
$ws.Cells.Item(37, 9).Value = "Primera"
$ws.Cells.Item(37, 10).Value = 50
$ws.Cells.Item(37, 11).Value = 9000
$ws.Cells.Item(37, 12).Value = 10000
$ws.Cells.Item(37, 13).Value = 9600
$ws.Cells.Item(37, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(37, 15).Value = "Región Metropolitana"
$ws.Cells.Item(37, 16).Value = 960
$ws.Cells.Item(37, 17).Value = 10
$ws.Cells.Item(37, 18).Value = "Hortaliza"
